$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '79.570.66'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").Value = '3.194.19'
$ws.Range("E3").Value = '  +5.24%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '207.12'
$ws.Range("E5").Value = '  +4.70%  '
$ws.Range("D6").Value = '635.44'
$ws.Range("E6").Value = '  +3.03%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.243'
$ws.Range("E8").Value = '  +19.07%  '
$ws.Range("D9").Value = '0.600'
$ws.Range("E9").Value = '  +9.90%  '
$ws.Range("D10").Value = '3.190.78'
$ws.Range("E10").Value = '  +5.19%  '
$ws.Range("D11").Value = '0.593'
$ws.Range("E11").Value = '  +37.13%  '
$ws.Range("D12").Value = '0.0000259'
$ws.Range("E12").Value = '  +35.74%  '
$ws.Range("E13").Value = '  +3.15%  '
$ws.Range("D14").Value = '5.40'
$ws.Range("E14").Value = '  +3.40%  '
$ws.Range("D15").Value = '3.779.15'
$ws.Range("E15").Value = '  +5.12%  '
$ws.Range("D16").Value = '31.97'
$ws.Range("E16").Value = '  +11.56%  '
$ws.Range("D17").Value = '79.366.34'
$ws.Range("E17").Value = '  +4.19%  '
$ws.Range("D18").Value = '3.189.73'
$ws.Range("E18").Value = '  +5.08%  '
$ws.Range("D19").Value = '14.54'
$ws.Range("E19").Value = '  +7.70%  '
$ws.Range("D20").Value = '9.36'
$ws.Range("E20").Value = '  +4.68%  '
$ws.Range("B21").Value = 'SuiNetwork'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D21").Value = '2.98'
$ws.Range("E21").Value = '  +28.30%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '440.16'
$ws.Range("E22").Value = '  +16.52%  '
$ws.Range("D23").Value = '5.17'
$ws.Range("E23").Value = '  +19.46%  '
$ws.Range("E24").Value = '  +12.10%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.358.04'
$ws.Range("E25").Value = '  +5.31%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '77.17'
$ws.Range("E26").Value = '  +6.36%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").Value = '10.82'
$ws.Range("E27").Value = '  +12.14%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0000123'
$ws.Range("E29").Value = '  +14.78%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '9.13'
$ws.Range("E30").Value = '  +10.99%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '1.53'
$ws.Range("E32").Value = '  +10.28%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '544.37'
$ws.Range("E33").Value = '  +11.04%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.159'
$ws.Range("E34").Value = '  +36.54%  '
$ws.Range("B35").Value = 'PancakeSwap'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D35").Value = '2.03'
$ws.Range("E35").Value = '  +6.32%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '23.07'
$ws.Range("E36").Value = '  +12.67%  '
$ws.Range("B37").Value = 'Cronos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D37").Value = '0.121'
$ws.Range("E37").Value = '  +16.67%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = '0.410'
$ws.Range("E39").Value = '  +8.13%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '163.94'
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '20.02'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '192.04'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").Value = '5.57'
$ws.Range("E44").Value = '  +10.47%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.82'
$ws.Range("E45").Value = '  +11.96%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.797'
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '2.66'
$ws.Range("E47").Value = '  +11.71%  '
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").Value = '1.33'
$ws.Range("E48").Value = '  +6.72%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '42.89'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").Value = '0.641'
$ws.Range("E50").Value = '  +7.39%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '25.70'
$ws.Range("E51").Value = '  +16.82%  '
